$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'63.815.23"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "'3.395.68"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'570.84"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "'161.82"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'3.396.97"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "'0.552"
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("D10").Value = "'7.29"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("E12").Value = "  -3.97%  "
$ws.Range("D13").Value = "'3.982.02"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "'26.93"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "'63.874.03"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "'3.346.60"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "'13.57"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").Value = "'375.30"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").Value = "'7.76"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").Value = "'70.23"
$ws.Range("E24").Value = "  -2.25%  "
$ws.Range("D25").Value = "'0.513"
$ws.Range("E25").Value = "  -5.25%  "
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("D27").Value = "'9.42"
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  -4.44%  "
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "'22.82"
$ws.Range("E33").Value = "  -1.09%  "
$ws.Range("D34").Value = "'7.07"
$ws.Range("E34").Value = "  +1.97%  "
$ws.Range("E35").Value = "  -4.59%  "
$ws.Range("D36").Value = "'159.68"
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").Value = "'0.858"
$ws.Range("E37").Value = "  +10.03%  "
$ws.Range("D38").Value = "'1.80"
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("D39").Value = "'0.0720"
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").Value = "'42.81"
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").Value = "'25.60"
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("D42").Value = "'6.45"
$ws.Range("E42").Value = "  -2.91%  "
$ws.Range("D43").Value = "'2.733.25"
$ws.Range("E43").Value = "  -5.90%  "
$ws.Range("D44").Value = "'26.07"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("D45").Value = "'4.37"
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("D46").Value = "'0.0305"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("E47").Value = "  +6.37%  "
$ws.Range("D48").Value = "'327.64"
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("D50").Value = "'6.27"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("E51").Value = "  -1.50%  "
